# ----------------------------------------------------------------------------
# Commit: [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# 1) "ODI Batting" sheet: the INNING_NUMBER (column B) cells that were
#    never actually scraped were left behind as blank placeholder cells.
#    Remove those now-empty cells outright (ClearContents -> no <c> at all),
#    leaving the populated INNING_NUMBER cells untouched.
# 2) Add a brand-new "ODI Batting Extra" sheet (after "ODI Bowling") holding
#    the newly-scraped extra batting fields (batting position, 4s, 6s,
#    percent of team runs, man-of-the-match) keyed by MATCH_CODE.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) Clean up the never-scraped INNING_NUMBER placeholder cells ----------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$emptyInningRows = @(3,4,5,6,7,8,9,10,11,12,15,16,17,18,20,22,23,25,27,29,30,39,41,42,43)
foreach ($r in $emptyInningRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# --- 2) Add the new "ODI Batting Extra" sheet, placed after "ODI Bowling" ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extraSheet = $wb.Worksheets.Add($null, $bowlingSheet)
$extraSheet.Name = "ODI Batting Extra"

# Header row: reuse the bold/bordered header style from "ODI Batting" column A
# (style index shared across the workbook's header rows), then set the text.
$battingSheet.Range("A1:F1").Copy()
$extraSheet.Range("A1:F1").PasteSpecial(-4122)

$extraSheet.Range("A1").Value = "MATCH_CODE"
$extraSheet.Range("B1").Value = "BATTING_POSITION"
$extraSheet.Range("C1").Value = "NUM_4"
$extraSheet.Range("D1").Value = "NUM_6"
$extraSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extraSheet.Range("F1").Value = "MAN_OF_MATCH"

# Data rows 2-43: A=MATCH_CODE, B=BATTING_POSITION, C=NUM_4, D=NUM_6,
# E=PERCENT_RUNS_OF_TOTAL, F=MAN_OF_MATCH.
# A leading "'" forces text storage for numeric-looking values (match codes,
# 4s/6s counts, percentages) so they round-trip as text, matching the
# scraped source data; an otherwise-empty field is written as a lone "'"
# so the cell exists (empty text) rather than being omitted.
$data = @(
    @("'3448", "'", "'", "'", "'", "NO"),
    @("'3449", "'", "'", "'", "'", "NO"),
    @("'3573", 10, "'", "'", "'", "NO"),
    @("'3576", 10, "'", "'", "'", "NO"),
    @("'3578", 10, "'", "'", "'", "NO"),
    @("'3596", "'", "'", "'", "'", "NO"),
    @("'3601", 10, "'", "'", "'", "NO"),
    @("'3713", 11, "'", "'", "'", "NO"),
    @("'3715", 9, "'", "'", "'", "NO"),
    @("'3719", "'", "'", "'", "'", "NO"),
    @("'3720", "'", "'", "'", "'", "NO"),
    @("'3723", 9, "'1", "'0", "'3.20%", "NO"),
    @("'3726", "'", "'", "'", "'", "NO"),
    @("'3734", 10, "'", "'", "'", "NO"),
    @("'3745", 10, "'", "'", "'", "NO"),
    @("'3747", 10, "'", "'", "'", "NO"),
    @("'3748", 10, "'", "'", "'", "NO"),
    @("'3753", 9, "'0", "'0", "'0.68%", "NO"),
    @("'3756", 9, "'", "'", "'", "NO"),
    @("'3767", 9, "'0", "'0", "'", "NO"),
    @("'3778", "'", "'", "'", "'", "NO"),
    @("'3793", 10, "'", "'", "'", "NO"),
    @("'3826", 9, "'1", "'0", "'5.28%", "NO"),
    @("'3827", 10, "'", "'", "'", "NO"),
    @("'3828", 9, "'0", "'0", "'1.81%", "NO"),
    @("'3865", 9, "'", "'", "'", "NO"),
    @("'3868", "'", "'", "'", "'", "NO"),
    @("'3870", "'", "'", "'", "'", "NO"),
    @("'3872", 9, "'", "'", "'", "NO"),
    @("'3883", 9, "'0", "'0", "'", "NO"),
    @("'3884", "'", "'", "'", "'", "NO"),
    @("'3886", 9, "'2", "'1", "'12.81%", "NO"),
    @("'3888", 9, "'1", "'0", "'2.03%", "NO"),
    @("'4026", 8, "'1", "'0", "'2.91%", "NO"),
    @("'4032", 9, "'2", "'0", "'3.78%", "NO"),
    @("'4036", 9, "'1", "'0", "'4.48%", "NO"),
    @("'4039", "'", "'", "'", "'", "NO"),
    @("'4085", "'", "'", "'", "'", "NO"),
    @("'4088", "'", "'", "'", "'", "NO"),
    @("'4089", 10, "'", "'", "'", "NO"),
    @("'4669", 8, "'", "'", "'", "NO"),
    @("'4676", 7, "'", "'", "'", "NO"),
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt 6; $c++) {
        $extraSheet.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

Write-Host "Added 'ODI Batting Extra' sheet and cleaned up blank INNING_NUMBER cells."
